$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "95.500.48"
$ws.Range("E2").Value = "  -1.82%  "

# Row 3
$ws.Range("D3").Value = "3.608.09"
$ws.Range("E3").Value = "  -2.43%  "

# Row 4
$ws.Range("D4").Value = "'2.67"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +23.46%  "

# Row 5
$ws.Range("E5").Value = "  +0.31%  "

# Row 6
$ws.Range("D6").Value = "'222.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.74%  "

# Row 7
$ws.Range("D7").Value = "'638.69"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.59%  "

# Row 8
$ws.Range("D8").Value = "'0.420"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.90%  "

# Row 9
$ws.Range("D9").Value = "'1.18"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.35%  "

# Row 10
$ws.Range("E10").Value = "  +0.13%  "

# Row 11
$ws.Range("D11").Value = "3.602.90"
$ws.Range("E11").Value = "  -2.52%  "

# Row 12
$ws.Range("D12").Value = "'48.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.69%  "

# Row 14
$ws.Range("D14").Value = "'0.0000289"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.38%  "

# Row 15
$ws.Range("D15").Value = "'6.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.82%  "

# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'26.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +39.45%  "

# Row 17
$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "4.280.46"
$ws.Range("E17").Value = "  -2.47%  "

# Row 18
$ws.Range("D18").Value = "95.387.78"
$ws.Range("E18").Value = "  -1.72%  "

# Row 19
$ws.Range("D19").Value = "'9.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.13%  "

# Row 20
$ws.Range("D20").Value = "'13.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.38%  "

# Row 21
$ws.Range("D21").Value = "3.601.73"
$ws.Range("E21").Value = "  -2.72%  "

# Row 22
$ws.Range("D22").Value = "'0.284"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +39.15%  "

# Row 23
$ws.Range("D23").Value = "'0.529"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.05%  "

# Row 24
$ws.Range("D24").Value = "'136.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +18.83%  "

# Row 25
$ws.Range("D25").Value = "'528.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.74%  "

# Row 26
$ws.Range("D26").Value = "'3.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.29%  "

# Row 27
$ws.Range("E27").Value = "  -8.78%  "

# Row 28
$ws.Range("E28").Value = "  -0.66%  "

# Row 29
$ws.Range("D29").Value = "3.776.44"
$ws.Range("E29").Value = "  -3.04%  "

# Row 30
$ws.Range("D30").Value = "'12.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.43%  "

# Row 31
$ws.Range("D31").Value = "'13.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.25%  "

# Row 32
$ws.Range("D32").Value = "'3.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.61%  "

# Row 33
$ws.Range("E33").Value = "  +0.13%  "

# Row 34
$ws.Range("D34").Value = "'0.639"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.18%  "

# Row 35
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "'1.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.95%  "

# Row 36
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'33.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.58%  "

# Row 37
$ws.Range("E37").Value = "  -3.16%  "

# Row 38
$ws.Range("E38").Value = "  +0.42%  "

# Row 40
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'7.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.66%  "

# Row 41
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "'587.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.11%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'8.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.03%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0526"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +16.35%  "

# Row 44
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.508"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.22%  "

# Row 45
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'1.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.66%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'41.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.19%  "

# Row 47
$ws.Range("D47").Value = "'0.159"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.48%  "

# Row 48
$ws.Range("D48").Value = "'1.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.24%  "

# Row 49
$ws.Range("D49").Value = "'9.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.54%  "

# Row 50
$ws.Range("D50").Value = "'234.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +13.93%  "

# Row 51
$ws.Range("E51").Value = "  -2.79%  "
